$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "customer_name" header (A1) and shift the remaining header cells
# (item_name ... location) one column to the left, carrying their formatting
# along with them. Row 2 (the first data row) is untouched.
$ws.Range("B1:K1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1:K1").Copy()
$ws.Range("A1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# The last header cell (K1) is now a stale duplicate of J1; remove it so the
# used range shrinks back down to column J.
$ws.Range("K1").Clear()

# "location" (now in column J) is wider than the old "remarks" header that
# used to live there, so the column grows a bit to fit it.
$ws.Columns.Item(10).ColumnWidth = 6.8

# Update the active cell selection as recorded in the edited workbook.
$ws.Range("G10").Select()
